$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells would otherwise be auto-converted from text to a number
# by Excel (e.g. "7.37" -> 7.37), unlike the original workbook where the Price
# column is stored as literal text. Force them to Text format first, cell by
# cell (a multi-cell/union Range only applies NumberFormat to its first area).
$textCells = @("D5", "D6", "D12", "D14", "D20", "D21", "D23", "D24", "D25", "D26", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D44", "D45", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) values.
$ws.Range("D2").Value = '51.327.62'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '2.978.03'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '382.17'
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").Value = '102.33'
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '0.0842'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '3.446.51'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").Value = '18.17'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").Value = '2.975.69'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("E17").Value = '  +5.14%  '
$ws.Range("D18").Value = '51.236.18'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '12.80'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").Value = '68.90'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '260.47'
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("D25").Value = '2.90'
$ws.Range("E25").Value = '  +6.39%  '
$ws.Range("D26").Value = '8.12'
$ws.Range("E26").Value = '  +13.08%  '
$ws.Range("E27").Value = '  +7.17%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '25.88'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = '9.82'
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = '34.35'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").Value = '50.84'
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").Value = '2.06'
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("D37").Value = '0.0453'
$ws.Range("E37").Value = '  +5.37%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").Value = '16.98'
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("E43").Value = '  -2.61%  '
$ws.Range("D44").Value = '122.88'
$ws.Range("D45").Value = '21.56'
$ws.Range("E45").Value = '  -2.93%  '
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("E47").Value = '  +3.03%  '
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").Value = '2.029.27'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").Value = '0.0332'
$ws.Range("E51").Value = '  +1.40%  '

# Restore the default (unstyled) cell style now that the values are committed
# as text, so only the text content differs from the original -- no leftover
# explicit "Text" number format on the cells.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}

